# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (date serial number)
$ws.Range("A2").Value = 45933

# Hourly prices 0h-1h ... 23h-24h
$ws.Range("B2").Value = 116.68
$ws.Range("C2").Value = 112.61
$ws.Range("D2").Value = 105.5
$ws.Range("E2").Value = 104.06
$ws.Range("F2").Value = 104.04
$ws.Range("G2").Value = 106.18
$ws.Range("H2").Value = 110.61
$ws.Range("I2").Value = 124.45
$ws.Range("J2").Value = 132.44
$ws.Range("K2").Value = 115.58
$ws.Range("L2").Value = 69.88
$ws.Range("M2").Value = 49.7
$ws.Range("N2").Value = 41.79
$ws.Range("O2").Value = 39.13
$ws.Range("P2").Value = 39.12
$ws.Range("Q2").Value = 38.23
$ws.Range("R2").Value = 44.1
$ws.Range("S2").Value = 58.39
$ws.Range("T2").Value = 94.91
$ws.Range("U2").Value = 119.98
$ws.Range("V2").Value = 138.76
$ws.Range("W2").Value = 134.48
$ws.Range("X2").Value = 120.81
$ws.Range("Y2").Value = 116.5
$ws.Range("Z2").Value = 93.25

# AA2 unchanged (Slot_4h_max = 20h-24h)
$ws.Range("AB2").Value = 127.64

# AC2 unchanged (Slot_2h_frist = 20h-22h)
$ws.Range("AD2").Value = 136.62

# Slot_2h_second changed from 6h-8h to 8h-10h
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 124.01

# Slot_min_price changed from 10h-18h to 10h-17h
$ws.Range("AG2").Value = "10h-17h"
